$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'96.115.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.97%  '
$ws.Range("D3").Value = "'3.584.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +8.96%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'239.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.57%  '
$ws.Range("D6").Value = "'636.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.11%  '
$ws.Range("E7").Value = '  +9.52%  '
$ws.Range("D8").Value = "'0.401"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.84%  '
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("D10").Value = "'1.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +11.94%  '
$ws.Range("D11").Value = "'3.582.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +9.05%  '
$ws.Range("D12").Value = "'43.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.98%  '
$ws.Range("E13").Value = '  +5.33%  '
$ws.Range("D14").Value = "'6.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +8.81%  '
$ws.Range("D15").Value = "'4.267.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +9.41%  '
$ws.Range("D16").Value = "'96.065.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.13%  '
$ws.Range("D17").Value = "'0.0000253"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.18%  '
$ws.Range("D18").Value = "'3.588.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +9.04%  '
$ws.Range("D19").Value = "'13.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +25.00%  '
$ws.Range("D20").Value = "'7.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.34%  '
$ws.Range("D21").Value = "'18.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.10%  '
$ws.Range("D22").Value = "'0.499"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +14.59%  '
$ws.Range("D23").Value = "'515.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.91%  '
$ws.Range("E24").Value = '  +2.09%  '
$ws.Range("E25").Value = '  +13.71%  '
$ws.Range("D26").Value = "'6.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +11.32%  '
$ws.Range("D27").Value = "'96.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +9.11%  '
$ws.Range("D28").Value = "'12.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.89%  '
$ws.Range("D29").Value = "'3.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +20.34%  '
$ws.Range("E30").Value = '  +7.01%  '
$ws.Range("D31").Value = "'11.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.21%  '
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("D33").Value = "'0.181"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.75%  '
$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.27%  '
$ws.Range("D35").Value = "'30.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +9.79%  '
$ws.Range("E36").Value = '  +9.63%  '
$ws.Range("D37").Value = "'579.72"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.46%  '
$ws.Range("D38").Value = "'7.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.01%  '
$ws.Range("E39").Value = '  +11.56%  '
$ws.Range("E40").Value = '  +4.94%  '
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").Value = "'0.927"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +10.09%  '
$ws.Range("D43").Value = "'0.0432"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.83%  '
$ws.Range("E44").Value = '  +6.71%  '
$ws.Range("D45").Value = "'23.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.23%  '
$ws.Range("E46").Value = '  +8.29%  '
$ws.Range("D47").Value = "'3.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.36%  '
$ws.Range("E48").Value = '  +6.05%  '
$ws.Range("D49").Value = "'53.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.45%  '
$ws.Range("D50").Value = "'8.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.53%  '
$ws.Range("D51").Value = "'3.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.47%  '
